$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the underlying engine allocates new shared-string table entries in the
# order distinct string values are first *written*, so the assignments below
# are deliberately ordered to reproduce the target sharedStrings.xml layout:
#   9=TSVData, 10=WebData, 11=Stats(LIKE), 12=Grants(LIKE), 13=Projects(LIKE),
#   14=Programs(LIKE), 15=Publications(LIKE)

# --- D2 / E2 stay the same text, but re-set first so they land right after the existing tail strings ---
$ws.Range("D2").Value = "TC07_INS_CancerType-EsophagealCancer_TSVData.xlsx"
$ws.Range("E2").Value = "TC07_INS_CancerType-EsophagealCancer_WebData.xlsx"

# --- Update C2: Stats query (now uses LIKE instead of IN) ---
$ws.Range("C2").Value = 'SELECT DISTINCT' + "`n" + `
'    COUNT(DISTINCT prg.program_id) AS "Programs",' + "`n" + `
'    COUNT(DISTINCT prj.project_id) AS "Projects",' + "`n" + `
'    COUNT(DISTINCT gnt.grant_id) AS "Grants",' + "`n" + `
'    COUNT(DISTINCT pub.pmid) AS "Publications"' + "`n" + `
'FROM ' + "`n" + `
'    df_program prg' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_project prj ON prg.program_id = prj."program.program_id"' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_grant gnt ON prj.project_id = gnt."project.project_id"' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_publication pub ON prj.project_id = pub."project.project_id"' + "`n" + `
'WHERE ' + "`n" + `
'    prg.cancer_type LIKE ' + "'" + '%Esophageal Cancer%' + "'" + ';'

# --- Update B4: Grants query (now uses LIKE, lower(grant_id) ORDER BY) ---
$ws.Range("B4").Value = 'SELECT DISTINCT' + "`n" + `
'    gnt.grant_id AS "Grant ID", ' + "`n" + `
'    prj.project_id AS "Project",' + "`n" + `
'    gnt.grant_title AS "Grant Title",' + "`n" + `
'    gnt.principal_investigators AS "Principal Investigators",' + "`n" + `
'    gnt.program_officers AS "Program Officers",' + "`n" + `
'    gnt.fiscal_year AS "Fiscal Year",' + "`n" + `
'    gnt.project_end_date AS "Project End Date"' + "`n" + `
'FROM ' + "`n" + `
'    df_grant gnt' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_project prj ON gnt."project.project_id" = prj.project_id' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_program prg ON prj."program.program_id" = prg.program_id' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_publication pub ON prj.project_id = pub."project.project_id"' + "`n" + `
'WHERE ' + "`n" + `
'    prg.cancer_type LIKE ' + "'" + '%Esophageal Cancer%' + "'" + "`n" + `
'ORDER BY ' + "`n" + `
'    lower(gnt.grant_id) ASC' + "`n" + `
'LIMIT 100;'

# --- Update B3: Projects query (now uses LIKE, lower(project_id) ORDER BY) ---
$ws.Range("B3").Value = 'SELECT DISTINCT' + "`n" + `
'    prj.project_id AS "Project ID", ' + "`n" + `
'    prj.project_title AS "Project Title",' + "`n" + `
'    prj.org_name AS "Organization",' + "`n" + `
'    prj.project_start_date AS "Project Start Date",' + "`n" + `
'    prj.project_end_date AS "Project End Date"' + "`n" + `
'FROM ' + "`n" + `
'    df_project prj' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_program prg ON prj."program.program_id" = prg.program_id' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_grant gnt ON prj.project_id = gnt."project.project_id"' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_publication pub ON prj.project_id = pub."project.project_id"' + "`n" + `
'WHERE ' + "`n" + `
'     prg.cancer_type LIKE ' + "'" + '%Esophageal Cancer%' + "'" + "`n" + `
'ORDER BY ' + "`n" + `
'    lower(prj.project_id) ASC' + "`n" + `
'LIMIT 100;'

# --- Update B2: Programs query (now uses LIKE, lower(program_name) ORDER BY, and CASE-based Data Location Details) ---
$ws.Range("B2").Value = 'SELECT DISTINCT ' + "`n" + `
'    prg.program_name AS "Program",' + "`n" + `
'    prg.website AS "Website",' + "`n" + `
'    prg.focus_area AS "Focus Area",' + "`n" + `
'    prg.cancer_type AS "Cancer Type",' + "`n" + `
'   CASE ' + "`n" + `
'        WHEN prg.data_link IS NOT NULL THEN prg.website       ' + "`n" + `
'        ELSE prg.data_link' + "`n" + `
'    END AS "Data Location Details"' + "`n" + `
'FROM ' + "`n" + `
'    df_program prg' + "`n" + `
'WHERE ' + "`n" + `
'     prg.cancer_type LIKE ' + "'" + '%Esophageal Cancer%' + "'" + "`n" + `
'ORDER BY ' + "`n" + `
'    lower(prg.program_name) ASC' + "`n" + `
'LIMIT 100;'

# --- Update B5: Publications query (now uses LIKE, lower(pmid) ORDER BY, new WHEN 3.0 branch) ---
$ws.Range("B5").Value = 'SELECT DISTINCT' + "`n" + `
'    pub.pmid AS "PubMed ID", ' + "`n" + `
'    pub.title AS "Title",' + "`n" + `
'    pub.authors AS "Authors",' + "`n" + `
'    pub.publication_date AS "Publication Date",' + "`n" + `
'    pub.cited_by AS "Cited By",' + "`n" + `
'    CASE ' + "`n" + `
'    WHEN pub.relative_citation_ratio = 0 THEN ' + "'0'" + "`n" + `
'    WHEN pub.relative_citation_ratio = 7.0 THEN ' + "'7'" + "`n" + `
'    WHEN pub.relative_citation_ratio = 2.0 THEN ' + "'2'" + "`n" + `
'    WHEN pub.relative_citation_ratio = 3.0 THEN ' + "'3'" + "`n" + `
'    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) ' + "`n" + `
'    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)' + "`n" + `
'END AS "Relative Citation Ratio"' + "`n" + `
'FROM ' + "`n" + `
'    df_publication pub' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_project prj ON pub."project.project_id" = prj.project_id' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_program prg ON prj."program.program_id" = prg.program_id' + "`n" + `
'LEFT JOIN ' + "`n" + `
'    df_grant gnt ON prj.project_id = gnt."project.project_id"' + "`n" + `
'WHERE ' + "`n" + `
'     prg.cancer_type LIKE ' + "'" + '%Esophageal Cancer%' + "'" + "`n" + `
'ORDER BY ' + "`n" + `
'    lower(pub.pmid) ASC' + "`n" + `
'LIMIT 100;'

# --- sheetView: move the selection to B2 (also resets the scrolled-in topLeftCell back to default) ---
$ws.Activate()
$ws.Range("B2").Select() | Out-Null

$wb.Save()
